# Weekly update: insert a new record row for "Poroto verde" at row 689,
# shifting all the following rows down by one (the last existing row
# ends up duplicated at the very bottom as row 778).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 689; Excel shifts rows 689:777 down
# to 690:778 and the sheet dimension grows to A1:R778 automatically.
$ws.Rows.Item(689).Insert()

# Match the date-cell style used by the surrounding rows (style index 2,
# i.e. the same number format as D690).
$ws.Range("D689").NumberFormat = $ws.Range("D690").NumberFormat

$ws.Range("A689").Value = 9
$ws.Range("B689").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C689").Value = "Metropolitana"
$ws.Range("D689").Value = 45077
$ws.Range("E689").Value = 13
$ws.Range("F689").Value = 100112031
$ws.Range("G689").Value = "Poroto verde"
$ws.Range("H689").Value = "Magnum"
$ws.Range("I689").Value = "Primera"
$ws.Range("J689").Value = 70
$ws.Range("K689").Value = 18000
$ws.Range("L689").Value = 20000
$ws.Range("M689").Value = 19000
$ws.Range("N689").Value = "`$/saco 25 kilos"
$ws.Range("O689").Value = "Región Metropolitana"
$ws.Range("P689").Value = 760
$ws.Range("Q689").Value = 25
$ws.Range("R689").Value = "Hortaliza"
